# Generate Report for Archive
#
# The localization status report is regenerated: every cell that held the
# "Ready for handoff" status (a single shared-string entry reused across the
# Overview sheet and each per-locale sheet) now reads "In Translation", and
# the Status column(s) are re-sized (narrower, since the new text is shorter
# than the old text) on every sheet that shows a Status column.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text -----------------------------------------------
# "Overview" shows the per-locale status in columns E (zh-cn) and F (de-de)
# for each of the 3 data rows; the locale sheets keep their own Status in
# column C. All of these cells share the same shared-string entry, so
# updating each cell's value re-points them all at the new text.
$overview.Range("E2:F4").Value = "In Translation"
$zhcn.Range("C2:C4").Value     = "In Translation"
$dede.Range("C2:C4").Value     = "In Translation"

# --- Re-size the Status column(s) to fit the new, shorter text ---------
# (ColumnWidth is quantized by the host to the nearest 1/6 character, so
# 12.5 is the input that lands on the nearest representable width to the
# regenerated report's column width.)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth     = 12.5
$dede.Columns.Item(3).ColumnWidth     = 12.5
